$wb = $excel.ActiveWorkbook

# --- Update the daily conversion summary text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.91 = 23655.78 pesos`n✅ 23655.78 pesos = 5.91 = 945.05 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 169.113
$wsTasas.Range("O10").Value = 4000.5
$wsTasas.Range("N12").Value = 4005
$wsTasas.Range("O12").Value = 160
